$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EER")

$ws.Range("A1").Value = 0.18909089850313116
$ws.Range("B1").Value = 0.1883158691412665
$ws.Range("C1").Value = 0.18950428968160471
$ws.Range("D1").Value = 0.18206111442472653
$ws.Range("E1").Value = 0.18689628297044972
$ws.Range("F1").Value = 0.18708011148677436

$ws.Range("A1:F1").Select() | Out-Null
